{"js": "// Correcao segundo orientacao do professor 1505\n//\n// 1) Move the \"_GoBack\" bookmark from the SSS_0022 paragraph (after\n//    \", exames\") to the very first paragraph (\"Clientes\"), wrapping\n//    that paragraph's text.\n// 2) Flip five requirement sentences from a negative phrasing\n//    (\"nao DEVE permitir ...\") to a positive one (\"DEVE bloquear ...\"):\n//    SSS_0002, SSS_0003, SSS_0004, SSS_0005, SSS_0008.\n\nconst body = context.document.body;\n\n// --- 1) Relocate the \"_GoBack\" bookmark ------------------------------\n// Remove it from wherever it currently is (if present) ...\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// ... and re-insert it around the text of the first paragraph (\"Clientes\").\nconst clientesResults = body.search(\"Clientes\", { matchCase: true });\nclientesResults.load(\"items\");\nawait context.sync();\nif (clientesResults.items.length > 0) {\n  clientesResults.items[0].insertBookmark(\"_GoBack\");\n}\n\n// --- 2) Flip the five \"nao DEVE permitir\" sentences to \"DEVE bloquear\" -\nconst replacements = [\n  [\n    \"SSS_0002 \u2013 O Sistema n\u00e3o DEVE permitir o cadastro de clientes menores de 18 anos.\",\n    \"SSS_0002 \u2013 O Sistema DEVE bloquear o cadastro de clientes menores de 18 anos.\"\n  ],\n  [\n    \"SSS_0003 \u2013 O Sistema n\u00e3o DEVE permitir o cadastro de um cliente sem animal.\",\n    \"SSS_0003 \u2013 O Sistema DEVE bloquear o cadastro de um cliente sem animal.\"\n  ],\n  [\n    \"SSS_0004 \u2013 O Sistema n\u00e3o DEVE permitir que o mesmo cliente seja cadastrado duas vezes.\",\n    \"SSS_0004 \u2013 O Sistema DEVE bloquear que o mesmo cliente seja cadastrado duas vezes.\"\n  ],\n  [\n    \"SSS_0005 \u2013 O Sistema n\u00e3o DEVE permitir que um usu\u00e1rio seja cadastrado sem CPF.\",\n    \"SSS_0005 \u2013 O Sistema DEVE bloquear que um usu\u00e1rio seja cadastrado sem CPF.\"\n  ],\n  [\n    \"SSS_0008 \u2013 O Sistema n\u00e3o DEVE permitir que um animal seja cadastrado sem um propriet\u00e1rio.\",\n    \"SSS_0008 \u2013 O Sistema DEVE bloquear que um animal seja cadastrado sem um propriet\u00e1rio.\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Correcao segundo orientacao do professor 1505\n#\n# 1) \"_GoBack\" bookmark moves from the SSS_0022 paragraph (after \", exames\")\n#    to the very first paragraph (\"Clientes\"), wrapping that paragraph's text.\n# 2) Five requirement sentences change from a negative phrasing (\"nao DEVE\n#    permitir ...\") to a positive one (\"DEVE bloquear ...\"):\n#       SSS_0002, SSS_0003, SSS_0004, SSS_0005, SSS_0008\n\n$d = $word.ActiveDocument\n\n# --- 1) Relocate the \"_GoBack\" bookmark -----------------------------------\n$existing = $d.Bookmarks\nif ($existing.Exists(\"_GoBack\")) {\n    $existing.Item(\"_GoBack\").Delete()\n}\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Clientes\")\n$d.Bookmarks.Add(\"_GoBack\", $rng)\n\n# --- 2) Flip the five \"nao DEVE permitir\" sentences to \"DEVE bloquear\" ----\n$pairs = @(\n    @(\"SSS_0002 \u2013 O Sistema n\u00e3o DEVE permitir o cadastro de clientes menores de 18 anos.\",\n      \"SSS_0002 \u2013 O Sistema DEVE bloquear o cadastro de clientes menores de 18 anos.\"),\n    @(\"SSS_0003 \u2013 O Sistema n\u00e3o DEVE permitir o cadastro de um cliente sem animal.\",\n      \"SSS_0003 \u2013 O Sistema DEVE bloquear o cadastro de um cliente sem animal.\"),\n    @(\"SSS_0004 \u2013 O Sistema n\u00e3o DEVE permitir que o mesmo cliente seja cadastrado duas vezes.\",\n      \"SSS_0004 \u2013 O Sistema DEVE bloquear que o mesmo cliente seja cadastrado duas vezes.\"),\n    @(\"SSS_0005 \u2013 O Sistema n\u00e3o DEVE permitir que um usu\u00e1rio seja cadastrado sem CPF.\",\n      \"SSS_0005 \u2013 O Sistema DEVE bloquear que um usu\u00e1rio seja cadastrado sem CPF.\"),\n    @(\"8 \u2013 O Sistema n\u00e3o DEVE permitir que um animal seja cadastrado sem um propriet\u00e1rio.\",\n      \"8 \u2013 O Sistema DEVE bloquear que um animal seja cadastrado sem um propriet\u00e1rio.\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $searchRange = $d.Content\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
